# Arbeitszeit.xlsx edit: update "Arbeitszeit 2" schedule start date + hours,
# reduce remaining vacation-day allowance, and update a logged work end-time
# in April, then let Excel recalculate all dependent formulas.

$wb = $excel.ActiveWorkbook

# --- Voreinstellungen sheet ---
$wsVor = $wb.Worksheets.Item("Voreinstellungen")

# "Arbeitszeit 2 ab" start date moves from 2025-03-01 to 2025-06-30
$wsVor.Range("B13").Value = 44376

# Monday hours for "Arbeitszeit 2" schedule (column D, row 13)
$wsVor.Range("D13").Value = 3.125

# Remaining vacation day entitlement for the year
$wsVor.Range("C36").Value = 0

# --- April sheet ---
$wsApril = $wb.Worksheets.Item("April")

# Logged end-time for April 1st (row 4, column E)
$wsApril.Range("E4").Value = 0.63541666666666663

# Recalculate all formulas across the workbook so dependent cells
# (monthly totals, Jahresübersicht, vacation-day messages, etc.) refresh.
$excel.CalculateFull()

# --- View-state touch-ups to mirror the saved selections ---
$wsVor.Range("E14").Select()

$wsMaerz = $wb.Worksheets.Item("März")
$wsMaerz.Activate()
$excel.ActiveWindow.ScrollRow = 4
$wsMaerz.Range("E34").Select()

$wsApril.Activate()
$wsApril.Range("E6").Select()
